$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10 (1-indexed), shifting rows 10:100 down to 11:101.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new faturamento diario record.
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 18313.97
$ws.Cells.Item(10, 3).Value = 7
$ws.Cells.Item(10, 4).Value = 2025
$ws.Cells.Item(10, 5).Value = "07/2025"
